$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Status Processo" values between M2 and M3
$ws.Range("M2").Value = "FATURADO"
$ws.Range("M3").Value = "PENDENTE"

# Update the active selection to A3
$ws.Range("A3").Select()
